# Applies the "three-digit number divided by one-digit number" worksheet edit:
# updates the date line and all the division problems in the table.
$d = $word.ActiveDocument

# wdReplaceAll = 2; used for every uniquely-occurring old value so we
# do not have to worry about locating a specific run/cell.
function Replace-AllText($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
        $true, 1, $false, $new, 2) | Out-Null
}

Replace-AllText "2024-06-02 Sunday" "2024-06-03 Monday"
Replace-AllText "184÷7=" "765÷5="
Replace-AllText "828÷8=" "605÷7="
Replace-AllText "390÷4=" "893÷8="
Replace-AllText "972÷7=" "530÷8="
Replace-AllText "831÷9=" "602÷8="
Replace-AllText "176÷7=" "110÷9="
Replace-AllText "902÷9=" "888÷2="
Replace-AllText "407÷2=" "285÷9="
Replace-AllText "465÷5=" "871÷2="
Replace-AllText "579÷2=" "513÷7="
Replace-AllText "748÷8=" "832÷4="
Replace-AllText "770÷2=" "793÷3="
Replace-AllText "153÷8=" "303÷4="
Replace-AllText "131÷2=" "990÷8="
Replace-AllText "359÷8=" "769÷3="
Replace-AllText "492÷4=" "643÷9="
Replace-AllText "615÷2=" "991÷6="
Replace-AllText "894÷7=" "613÷7="
Replace-AllText "244÷8=" "100÷7="
Replace-AllText "500÷6=" "191÷5="
Replace-AllText "644÷8=" "370÷4="
Replace-AllText "359÷6=" "869÷4="
Replace-AllText "392÷9=" "995÷8="

# "450÷4=" occurs twice in the table (row 9 and row 17, column 5) and the two
# occurrences must become different values, so address each table cell
# directly and replace only the first match within that cell's range
# (wdReplaceOne = 1) rather than using the document-wide replace above.
$t = $d.Tables.Item(1)
$cell = $t.Cell(9, 5)
$cell.Range.Find.Execute("450÷4=", $true, $false, $false, $false, $false, `
    $true, 0, $false, "117÷2=", 1) | Out-Null
$cell = $t.Cell(17, 5)
$cell.Range.Find.Execute("450÷4=", $true, $false, $false, $false, $false, `
    $true, 0, $false, "567÷8=", 1) | Out-Null
